$d = $word.ActiveDocument

# Locate the paragraph that currently ends with "... Experimentator als Endboss."
$para = $d.Paragraphs(7)
$rng = $para.Range

# Move the end back before the paragraph mark so we insert inside the paragraph,
# then collapse to that point (i.e. right after "Endboss.").
$rng.MoveEnd(1, -1) | Out-Null
$rng.Collapse(0) | Out-Null
$insertPos = $rng.Start

# Insert the new sentence as a new run right after the existing text.
$newText = " Server nicht zerstören, da sich Matrix schon auf mehrere ausgebreitet hat. Es wäre ein zu großer Schaden."
$rng.InsertAfter($newText)

# Give the newly inserted run the same German language formatting used elsewhere
# in the document (re-use the same Range object that did the insertion so the
# formatting lands on the freshly created run).
$rng.LanguageID = "de-DE"
$newTextEnd = $rng.End

# Append a temporary placeholder character. This keeps the spot where the
# bookmark needs to go from being the very last character position of the
# paragraph (collapsed bookmarks placed exactly there resolve incorrectly),
# then it is removed again once the bookmark has been created.
$rng.InsertAfter("X")

# Move the "_GoBack" bookmark from the end of the document (after "Anzahl NPCs: 11")
# to right after the text we just inserted.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bookmarkRange = $d.Range($newTextEnd, $newTextEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary placeholder character again.
$placeholderRange = $d.Range($newTextEnd, $newTextEnd + 1)
$placeholderRange.Delete()
